# New crime data collected - weekly refresh for cs-en-us-122pct.xlsx
# Updates the "Volume/Number" + reporting-week header text and refreshes
# the weekly/28-day/YTD/2yr crime-complaint figures (and their computed
# % changes) in the data table (rows 14-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: Volume 30 Number 32 -> 33; week of 8/7-8/13 -> 8/14-8/20
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  33"
$ws.Range("C9").Value = "Report Covering the Week  8/14/2023  Through  8/20/2023"

# ---------------------------------------------------------------------
# Helper: some cells flip between the "number" style and the "N/A text"
# style (e.g. a count drops to/from 0 and the %-chg column shows the
# "***.*" placeholder instead of a computed number). Copying formats
# from an untouched donor cell in the same row keeps font/alignment
# identical while giving us the right style bucket for the new value.
# ---------------------------------------------------------------------

# Row 15 - Rape: 2023 WTD count now 0 (text "0"), % chg -> "***.*"
$ws.Range("A15").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("D15").Value = "0"
$ws.Range("E15").Value = "***.*"
$ws.Range("M15").Value2 = -53.846153846153

# Row 16 - Robbery
$ws.Range("C16").Value2 = 2
$ws.Range("A16").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").Value = "0"
$ws.Range("E16").Value = "***.*"
$ws.Range("F16").Value2 = 6
$ws.Range("H16").Value2 = 100
$ws.Range("I16").Value2 = 37
$ws.Range("K16").Value2 = -2.631578947368
$ws.Range("L16").Value2 = 37.037037037037
$ws.Range("M16").Value2 = -43.076923076923
$ws.Range("N16").Value2 = -78.857142857142

# Row 17 - Fel. Assault
$ws.Range("D17").Value2 = 3
$ws.Range("E17").Value2 = 33.333333333333
$ws.Range("F17").Value2 = 15
$ws.Range("H17").Value2 = -11.764705882352
$ws.Range("I17").Value2 = 120
$ws.Range("J17").Value2 = 73
$ws.Range("K17").Value2 = 64.383561643835
$ws.Range("L17").Value2 = 55.844155844155
$ws.Range("M17").Value2 = 34.831460674157
$ws.Range("N17").Value2 = -38.461538461538

# Row 18 - Burglary
$ws.Range("D18").Value2 = 2
$ws.Range("E18").Value2 = 0
$ws.Range("F18").Value2 = 8
$ws.Range("G18").Value2 = 4
$ws.Range("H18").Value2 = 100
$ws.Range("I18").Value2 = 66
$ws.Range("J18").Value2 = 42
$ws.Range("K18").Value2 = 57.142857142857
$ws.Range("L18").Value2 = 34.693877551020
$ws.Range("M18").Value2 = -48.031496062992
$ws.Range("N18").Value2 = -92.271662763466

# Row 19 - Gr. Larceny
$ws.Range("D19").Value2 = 4
$ws.Range("E19").Value2 = 125
$ws.Range("F19").Value2 = 35
$ws.Range("G19").Value2 = 28
$ws.Range("H19").Value2 = 25
$ws.Range("I19").Value2 = 299
$ws.Range("J19").Value2 = 210
$ws.Range("K19").Value2 = 42.380952380952
$ws.Range("L19").Value2 = 41.037735849056
$ws.Range("M19").Value2 = 15.444015444015
$ws.Range("N19").Value2 = -43.796992481203

# Row 20 - G.L.A.
$ws.Range("C20").Value2 = 3
$ws.Range("D20").Value2 = 2
$ws.Range("E20").Value2 = 50
$ws.Range("F20").Value2 = 15
$ws.Range("G20").Value2 = 13
$ws.Range("H20").Value2 = 15.384615384615
$ws.Range("I20").Value2 = 75
$ws.Range("J20").Value2 = 84
$ws.Range("K20").Value2 = -10.714285714285
$ws.Range("L20").Value2 = 134.375
$ws.Range("M20").Value2 = 7.142857142857
$ws.Range("N20").Value2 = -95.810055865921

# Row 21 - TOTAL
$ws.Range("C21").Value2 = 20
$ws.Range("D21").Value2 = 11
$ws.Range("E21").Value2 = 81.818181818181
$ws.Range("F21").Value2 = 79
$ws.Range("G21").Value2 = 67
$ws.Range("H21").Value2 = 17.910447761194
$ws.Range("I21").Value2 = 604
$ws.Range("J21").Value2 = 454
$ws.Range("K21").Value2 = 33.039647577092
$ws.Range("L21").Value2 = 49.504950495049
$ws.Range("M21").Value2 = -3.049759229534
$ws.Range("N21").Value2 = -83.057503506311

# Row 23 - Housing
$ws.Range("C23").Value2 = 1
$ws.Range("F23").Value2 = 6
$ws.Range("H23").Value2 = 500
$ws.Range("I23").Value2 = 27
$ws.Range("K23").Value2 = 107.692307692308
$ws.Range("L23").Value2 = -3.571428571428
$ws.Range("M23").Value2 = 107.692307692308

# Row 24 - Petit Larceny
$ws.Range("C24").Value2 = 17
$ws.Range("D24").Value2 = 22
$ws.Range("E24").Value2 = -22.727272727272
$ws.Range("G24").Value2 = 89
$ws.Range("H24").Value2 = 0
$ws.Range("I24").Value2 = 717
$ws.Range("J24").Value2 = 496
$ws.Range("K24").Value2 = 44.556451612903
$ws.Range("L24").Value2 = 142.22972972973
$ws.Range("M24").Value2 = -35.579514824797

# Row 25 - Misd. Assault
$ws.Range("D25").Value2 = 5
$ws.Range("E25").Value2 = 60
$ws.Range("F25").Value2 = 34
$ws.Range("G25").Value2 = 25
$ws.Range("H25").Value2 = 36
$ws.Range("I25").Value2 = 209
$ws.Range("J25").Value2 = 221
$ws.Range("K25").Value2 = -5.429864253393
$ws.Range("L25").Value2 = 36.601307189542
$ws.Range("M25").Value2 = -41.782729805013

# Row 26 - UCR Rape*: 2023 WTD now 1 (was "0" text), 2022 WTD now "0" text
$ws.Range("G26").Copy() | Out-Null
$ws.Range("C26").PasteSpecial(-4122) | Out-Null
$ws.Range("F26").PasteSpecial(-4122) | Out-Null
$ws.Range("C26").Value2 = 1
$ws.Range("F26").Value2 = 1
$ws.Range("A26").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null
$ws.Range("D26").Value = "0"
$ws.Range("E26").Value = "***.*"
$ws.Range("H26").Value2 = -66.666666666666
$ws.Range("I26").Value2 = 10
$ws.Range("K26").Value2 = -28.571428571428
$ws.Range("L26").Value2 = 0

# Row 27 - Other Sex Crimes: 2022 WTD now 0 (text), % chg "***.*"
$ws.Range("A27").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("D27").Value = "0"
$ws.Range("E27").Value = "***.*"

# Row 30 - Hate Crimes: 28-day 2022 now 0 (text), % chg "***.*"
$ws.Range("A30").Copy() | Out-Null
$ws.Range("G30").PasteSpecial(-4122) | Out-Null
$ws.Range("H30").PasteSpecial(-4122) | Out-Null
$ws.Range("G30").Value = "0"
$ws.Range("H30").Value = "***.*"
$ws.Range("L30").Value2 = 50
